$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D price cells to remain text (matches source data which is
# inline-string, not numeric) while writing the updated values, then restore
# the default "Normal" style so no stray number formatting is left behind.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "37.364.40"
$ws.Range("E2").Value = "  -1.21%  "
$ws.Range("D3").Value = "2.051.76"
$ws.Range("E3").Value = "  -1.27%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "230.64"
$ws.Range("E5").Value = "  -1.32%  "
$ws.Range("D6").Value = "0.622"
$ws.Range("E6").Value = "  -0.56%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "57.18"
$ws.Range("E8").Value = "  -3.48%  "
$ws.Range("E9").Value = "  -2.34%  "
$ws.Range("E10").Value = "  -2.42%  "
$ws.Range("E11").Value = "  +1.42%  "
$ws.Range("E12").Value = "  -0.33%  "
$ws.Range("D13").Value = "2.354.17"
$ws.Range("E13").Value = "  -1.26%  "
$ws.Range("D14").Value = "20.61"
$ws.Range("E14").Value = "  -2.95%  "
$ws.Range("D15").Value = "0.757"
$ws.Range("E15").Value = "  -2.36%  "
$ws.Range("D16").Value = "5.26"
$ws.Range("E16").Value = "  -1.99%  "
$ws.Range("D17").Value = "2.056.98"
$ws.Range("E17").Value = "  -1.13%  "
$ws.Range("D18").Value = "37.373.81"
$ws.Range("E18").Value = "  -0.90%  "
$ws.Range("D19").Value = "5.99"
$ws.Range("E19").Value = "  -2.65%  "
$ws.Range("D20").Value = "69.81"
$ws.Range("E20").Value = "  -2.49%  "
$ws.Range("E21").Value = "  -3.47%  "
$ws.Range("D22").Value = "226.57"
$ws.Range("E22").Value = "  -0.76%  "
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("D24").Value = "2.39"
$ws.Range("E24").Value = "  +1.31%  "
$ws.Range("E25").Value = "  -3.33%  "
$ws.Range("E26").Value = "  +3.63%  "
$ws.Range("D27").Value = "168.83"
$ws.Range("E27").Value = "  -1.17%  "
$ws.Range("E28").Value = "  -3.08%  "
$ws.Range("D29").Value = "19.20"
$ws.Range("E29").Value = "  -1.63%  "
$ws.Range("D30").Value = "1.36"
$ws.Range("E30").Value = "  -4.52%  "
$ws.Range("E31").Value = "  +0.42%  "
$ws.Range("E32").Value = "  -3.79%  "
$ws.Range("D33").Value = "0.0626"
$ws.Range("E33").Value = "  -1.19%  "
$ws.Range("D34").Value = "4.57"
$ws.Range("D35").Value = "2.48"
$ws.Range("E35").Value = "  -0.33%  "
$ws.Range("E36").Value = "  -0.34%  "
$ws.Range("E37").Value = "  -3.63%  "
$ws.Range("E38").Value = "  -0.16%  "
$ws.Range("E39").Value = "  -2.17%  "
$ws.Range("E40").Value = "  +5.10%  "
$ws.Range("D41").Value = "98.10"
$ws.Range("E41").Value = "  -1.26%  "
$ws.Range("D42").Value = "0.0955"
$ws.Range("E42").Value = "  -2.07%  "
$ws.Range("E43").Value = "  +0.62%  "
$ws.Range("D44").Value = "1.482.49"
$ws.Range("E44").Value = "  +2.60%  "
$ws.Range("E45").Value = "  +3.38%  "
$ws.Range("D46").Value = "16.63"
$ws.Range("E46").Value = "  +0.22%  "
$ws.Range("D47").Value = "4.04"
$ws.Range("E47").Value = "  -3.04%  "
$ws.Range("E48").Value = "  -3.18%  "
$ws.Range("D49").Value = "7.24"
$ws.Range("E49").Value = "  -2.36%  "
$ws.Range("E50").Value = "  -3.02%  "
$ws.Range("D51").Value = "2.239.95"
$ws.Range("E51").Value = "  -1.28%  "

$ws.Range("D2:D51").Style = "Normal"
